# Scheduled-runner market data refresh for Sheets/Atomos_Profits.xlsx
# Updates currentAveragePrice* / Leve*Price / LeveProfit* columns (H:N) per leve row
# across all 8 crafting-job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()
$ws.Range("H116").Value = 3891.4583
$ws.Range("I116").Value = 3134.9092
$ws.Range("K116").Value = 3134.9092
$ws.Range("M116").Value = 307.0907999999999
$ws.Range("H132").Value = 5410727.5
$ws.Range("I132").Value = 6671980.5
$ws.Range("J132").Value = 5358.857
$ws.Range("K132").Value = 20015941.5
$ws.Range("L132").Value = 16076.571
$ws.Range("M132").Value = -20013411.5
$ws.Range("N132").Value = -21136.571

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 15627365
$ws.Range("I2").Value = 25001294
$ws.Range("J2").Value = 4150
$ws.Range("K2").Value = 25001294
$ws.Range("L2").Value = 4150
$ws.Range("M2").Value = -25001181
$ws.Range("N2").Value = -4376
$ws.Range("H5").Value = 245.53847
$ws.Range("I5").Value = 196
$ws.Range("J5").Value = 276.5
$ws.Range("K5").Value = 196
$ws.Range("L5").Value = 276.5
$ws.Range("M5").Value = -84
$ws.Range("N5").Value = -500.5
$ws.Range("H45").Value = 1800.6897
$ws.Range("I45").Value = 1090.4166
$ws.Range("K45").Value = 1090.4166
$ws.Range("M45").Value = -713.4166
$ws.Range("H52").Value = 39420
$ws.Range("J52").Value = 39420
$ws.Range("L52").Value = 39420
$ws.Range("N52").Value = -40056
$ws.Range("H74").Value = 1004.6667
$ws.Range("I74").Value = 1000
$ws.Range("J74").Value = 1007
$ws.Range("K74").Value = 1000
$ws.Range("L74").Value = 1007
$ws.Range("M74").Value = -126
$ws.Range("N74").Value = -2755
$ws.Range("H77").Value = 1004.6667
$ws.Range("I77").Value = 1000
$ws.Range("J77").Value = 1007
$ws.Range("K77").Value = 5000
$ws.Range("L77").Value = 5035
$ws.Range("M77").Value = -632
$ws.Range("N77").Value = -13771
$ws.Range("H102").Value = 2589.2307
$ws.Range("I102").Value = 2596.6667
$ws.Range("J102").Value = 2500
$ws.Range("K102").Value = 2596.6667
$ws.Range("L102").Value = 2500
$ws.Range("M102").Value = -974.6667000000002
$ws.Range("N102").Value = -5744
$ws.Range("H110").Value = 1520.5
$ws.Range("I110").Value = 715.53845
$ws.Range("J110").Value = 2471.818
$ws.Range("K110").Value = 715.53845
$ws.Range("L110").Value = 2471.818
$ws.Range("M110").Value = 1329.46155
$ws.Range("N110").Value = -6561.818
$ws.Range("H116").Value = 15627365
$ws.Range("I116").Value = 25001294
$ws.Range("J116").Value = 4150
$ws.Range("K116").Value = 25001294
$ws.Range("L116").Value = 4150
$ws.Range("M116").Value = -24999000
$ws.Range("N116").Value = -8738
$ws.Range("H122").Value = 3302.7896
$ws.Range("I122").Value = 2211.7693
$ws.Range("K122").Value = 6635.3079
$ws.Range("M122").Value = -4185.3079
$ws.Range("H132").Value = 23813850
$ws.Range("I132").Value = 34486896
$ws.Range("J132").Value = 4748.154
$ws.Range("K132").Value = 103460688
$ws.Range("L132").Value = 14244.462
$ws.Range("M132").Value = -103458158
$ws.Range("N132").Value = -19304.462
$ws.Range("H134").Value = 97414.5
$ws.Range("J134").Value = 97414.5
$ws.Range("L134").Value = 97414.5
$ws.Range("N134").Value = -107554.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 15627365
$ws.Range("I3").Value = 25001294
$ws.Range("J3").Value = 4150
$ws.Range("K3").Value = 25001294
$ws.Range("L3").Value = 4150
$ws.Range("M3").Value = -25001180
$ws.Range("N3").Value = -4378
$ws.Range("H4").Value = 245.53847
$ws.Range("I4").Value = 196
$ws.Range("J4").Value = 276.5
$ws.Range("K4").Value = 196
$ws.Range("L4").Value = 276.5
$ws.Range("M4").Value = -81
$ws.Range("N4").Value = -506.5
$ws.Range("H94").Value = 635.2222
$ws.Range("I94").Value = 616.7143
$ws.Range("K94").Value = 616.7143
$ws.Range("M94").Value = -165.7143
$ws.Range("H99").Value = 1799.7368
$ws.Range("I99").Value = 1256.7858
$ws.Range("J99").Value = 3320
$ws.Range("K99").Value = 1256.7858
$ws.Range("L99").Value = 3320
$ws.Range("M99").Value = 241.2141999999999
$ws.Range("N99").Value = -6316
$ws.Range("H105").Value = 1676.375
$ws.Range("I105").Value = 1472.7778
$ws.Range("J105").Value = 2287.1667
$ws.Range("K105").Value = 1472.7778
$ws.Range("L105").Value = 2287.1667
$ws.Range("M105").Value = 274.2221999999999
$ws.Range("N105").Value = -5781.1667
$ws.Range("H107").Value = 2483.6956
$ws.Range("I107").Value = 2131.0833
$ws.Range("K107").Value = 2131.0833
$ws.Range("M107").Value = -211.0832999999998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 191.125
$ws.Range("I7").Value = 35
$ws.Range("J7").Value = 347.25
$ws.Range("K7").Value = 35
$ws.Range("L7").Value = 347.25
$ws.Range("M7").Value = 78
$ws.Range("N7").Value = -573.25
$ws.Range("H22").Value = 2328.7144
$ws.Range("I22").Value = 700.5
$ws.Range("J22").Value = 2980
$ws.Range("K22").Value = 700.5
$ws.Range("L22").Value = 2980
$ws.Range("M22").Value = -350.5
$ws.Range("N22").Value = -3680
$ws.Range("H99").Value = 2219.4614
$ws.Range("I99").Value = 1172
$ws.Range("J99").Value = 3441.5
$ws.Range("K99").Value = 1172
$ws.Range("L99").Value = 3441.5
$ws.Range("M99").Value = 326
$ws.Range("N99").Value = -6437.5
$ws.Range("H105").Value = 2886.8572
$ws.Range("I105").Value = 1827
$ws.Range("J105").Value = 4300
$ws.Range("K105").Value = 1827
$ws.Range("L105").Value = 4300
$ws.Range("M105").Value = -80
$ws.Range("N105").Value = -7794
$ws.Range("H126").Value = 2219.4614
$ws.Range("I126").Value = 1172
$ws.Range("J126").Value = 3441.5
$ws.Range("K126").Value = 3516
$ws.Range("L126").Value = 10324.5
$ws.Range("M126").Value = -1046
$ws.Range("N126").Value = -15264.5
$ws.Range("H134").Value = 1751.0952
$ws.Range("I134").Value = 951.64703
$ws.Range("J134").Value = 5148.75
$ws.Range("K134").Value = 2854.94109
$ws.Range("L134").Value = 15446.25
$ws.Range("M134").Value = -319.9410899999998
$ws.Range("N134").Value = -20516.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 5991.7144
$ws.Range("I132").Value = 848
$ws.Range("J132").Value = 9849.5
$ws.Range("K132").Value = 7632
$ws.Range("L132").Value = 88645.5
$ws.Range("M132").Value = -5102
$ws.Range("N132").Value = -93705.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1903
$ws.Range("I97").Value = 1053.75
$ws.Range("J97").Value = 5300
$ws.Range("K97").Value = 1053.75
$ws.Range("L97").Value = 5300
$ws.Range("M97").Value = -557.75
$ws.Range("N97").Value = -6292
$ws.Range("H102").Value = 2383.6572
$ws.Range("I102").Value = 1893.12
$ws.Range("J102").Value = 3610
$ws.Range("K102").Value = 1893.12
$ws.Range("L102").Value = 3610
$ws.Range("M102").Value = -271.1199999999999
$ws.Range("N102").Value = -6854
$ws.Range("H107").Value = 883.4231
$ws.Range("I107").Value = 363.47058
$ws.Range("K107").Value = 363.47058
$ws.Range("M107").Value = 1556.52942
$ws.Range("H126").Value = 2721.4348
$ws.Range("I126").Value = 1864.6666
$ws.Range("K126").Value = 5593.9998
$ws.Range("M126").Value = -3123.9998
$ws.Range("H132").Value = 3159.3142
$ws.Range("I132").Value = 2503.762
$ws.Range("J132").Value = 4142.643
$ws.Range("K132").Value = 7511.286
$ws.Range("L132").Value = 12427.929
$ws.Range("M132").Value = -4981.286
$ws.Range("N132").Value = -17487.929

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2746.5
$ws.Range("I40").Value = 1687
$ws.Range("J40").Value = 3099.6667
$ws.Range("K40").Value = 1687
$ws.Range("L40").Value = 3099.6667
$ws.Range("M40").Value = -1551
$ws.Range("N40").Value = -3371.6667
$ws.Range("H61").Value = 43479944
$ws.Range("I61").Value = 50001280
$ws.Range("J61").Value = 4366.6665
$ws.Range("K61").Value = 50001280
$ws.Range("L61").Value = 4366.6665
$ws.Range("M61").Value = -50001078
$ws.Range("N61").Value = -4770.6665
$ws.Range("H82").Value = 3131.5
$ws.Range("J82").Value = 4100
$ws.Range("L82").Value = 4100
$ws.Range("N82").Value = -4822
$ws.Range("H85").Value = 3131.5
$ws.Range("J85").Value = 4100
$ws.Range("L85").Value = 4100
$ws.Range("N85").Value = -6596
$ws.Range("H100").Value = 2285.4546
$ws.Range("I100").Value = 1350
$ws.Range("J100").Value = 2820
$ws.Range("K100").Value = 1350
$ws.Range("L100").Value = 2820
$ws.Range("M100").Value = -809
$ws.Range("N100").Value = -3902
$ws.Range("H113").Value = 43479944
$ws.Range("I113").Value = 50001280
$ws.Range("J113").Value = 4366.6665
$ws.Range("K113").Value = 50001280
$ws.Range("L113").Value = 4366.6665
$ws.Range("M113").Value = -49999110
$ws.Range("N113").Value = -8706.666499999999
$ws.Range("H122").Value = 2983.7297
$ws.Range("I122").Value = 2592.5557
$ws.Range("K122").Value = 7777.6671
$ws.Range("M122").Value = -5327.6671
$ws.Range("H127").Value = 31694
$ws.Range("J127").Value = 31694
$ws.Range("L127").Value = 31694
$ws.Range("N127").Value = -41614
$ws.Range("H132").Value = 3469.5386
$ws.Range("I132").Value = 2015
$ws.Range("J132").Value = 4005.4211
$ws.Range("K132").Value = 6045
$ws.Range("L132").Value = 12016.2633
$ws.Range("M132").Value = -3515
$ws.Range("N132").Value = -17076.2633

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 18114.834
$ws.Range("I96").Value = 1672.25
$ws.Range("J96").Value = 51000
$ws.Range("K96").Value = 1672.25
$ws.Range("L96").Value = 51000
$ws.Range("M96").Value = -299.25
$ws.Range("N96").Value = -53746
$ws.Range("H132").Value = 249596
$ws.Range("I132").Value = 336214.6
$ws.Range("J132").Value = 13363.454
$ws.Range("K132").Value = 1008643.8
$ws.Range("L132").Value = 40090.362
$ws.Range("M132").Value = -1006113.8
$ws.Range("N132").Value = -45150.362
